$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "JD_003"
$ws.Range("B4").Value = "Cyber Security Engineer"
$ws.Range("C4").Value = "efesafef"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = "Hybrid"
$ws.Range("G4").Value = "Pune, Maharashtra, India"
